$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# Sales Order Identifier columns (S) and Previous Doc columns (AX) -
# numeric-looking text, must stay text (not auto-converted to a number)
Set-TextValue "S2"  "7601959556"
Set-TextValue "AX2" "7601959556"
Set-TextValue "S3"  "7601959556"
Set-TextValue "AX3" "7601959556"

# Financial Close Calendar Date columns (T) - DD/MM/YYYY text, not a
# recognizable date in this locale so it is stored as plain text already.
$ws.Range("T2").Value = "18/08/2016"
$ws.Range("T3").Value = "31/08/2016"

$excel.CutCopyMode = 0
